$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the taxon-record data between row 3 and row 4 (A, B, E, F, G, H),
# and move the "Publik kommentar" (AC) comment from row 4 to row 3.

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$e4 = $ws.Range("E4").Value2
$f4 = $ws.Range("F4").Value2
$g4 = $ws.Range("G4").Value2
$h4 = $ws.Range("H4").Value2
$ac4 = $ws.Range("AC4").Value2

$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("E3").Value = $e4
$ws.Range("F3").Value = $f4
$ws.Range("G3").Value = $g4
$ws.Range("H3").Value = $h4
$ws.Range("AC3").Value = $ac4

$ws.Range("A4").Value = $a3
$ws.Range("B4").Value = $b3
$ws.Range("E4").Value = $e3
$ws.Range("F4").Value = $f3
$ws.Range("G4").Value = $g3
$ws.Range("H4").Value = $h3
$ws.Range("AC4").ClearContents()
